$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -21.853
    18 = -22.095
    20 = -20.623
    27 = -21.918
    35 = -20.123
    69 = -21.52
    76 = -20.047
    78 = -19.854
    82 = -21.962
    83 = -21.891
    93 = -21.508
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}

$wb.Save()
